$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated figures for the 2022-05-16 data refresh of the
# 'Fonds de solidarite - volet 1' dataset (nombre_aides / nombre_entreprises / montant_total).
$ws.Range("C64").Value = 5203
$ws.Range("E64").Value = 20360065
$ws.Range("C66").Value = 768
$ws.Range("E66").Value = 9929275
$ws.Range("C70").Value = 15726
$ws.Range("E70").Value = 24658731
$ws.Range("C73").Value = 2456
$ws.Range("E73").Value = 7383072
$ws.Range("C74").Value = 948
$ws.Range("E74").Value = 4258463
$ws.Range("C75").Value = 403
$ws.Range("E75").Value = 2849572
$ws.Range("C76").Value = 128
$ws.Range("E76").Value = 2483190
$ws.Range("C91").Value = 151098
$ws.Range("E91").Value = 482084331
$ws.Range("C92").Value = 409006
$ws.Range("E92").Value = 1593616902
$ws.Range("C93").Value = 209490
$ws.Range("E93").Value = 1307691204
$ws.Range("C94").Value = 94144
$ws.Range("E94").Value = 915770565
$ws.Range("C95").Value = 50722
$ws.Range("E95").Value = 930422875
$ws.Range("C98").Value = 810
$ws.Range("E98").Value = 117791167
$ws.Range("C101").Value = 179
$ws.Range("E101").Value = 32052522
$ws.Range("C104").Value = 135216
$ws.Range("E104").Value = 272104589
$ws.Range("C105").Value = 8170
$ws.Range("E105").Value = 16872184
$ws.Range("C107").Value = 6390
$ws.Range("E107").Value = 21955635
$ws.Range("C114").Value = 3798
$ws.Range("E114").Value = 9098224
$ws.Range("C116").Value = 4558
$ws.Range("E116").Value = 20484785
$ws.Range("C117").Value = 1913
$ws.Range("E117").Value = 12352831
$ws.Range("C118").Value = 976
$ws.Range("E118").Value = 11793478
$ws.Range("C122").Value = 8488
$ws.Range("E122").Value = 12672414
$ws.Range("C132").Value = 30288
$ws.Range("E132").Value = 174222945
$ws.Range("C144").Value = 24413
$ws.Range("E144").Value = 201797181
$ws.Range("C173").Value = 96858
$ws.Range("E173").Value = 327924730
$ws.Range("C184").Value = 68734
$ws.Range("D184").Value = 13881
$ws.Range("E184").Value = 134170601
